# Informe Diario - actualizacion del informe diario
#
# 1. Move the hidden "_GoBack" bookmark away from the end of the
#    "Soluciones:  Pong ya inicializa" paragraph (Dia 26 section).
# 2. Append new text to the "Problemas:" paragraph in the "Diario Dia 28"
#    section: " Poner imagen en los botones donde estan las bombas."
# 3. Re-create the "_GoBack" bookmark, now collapsed right after that new
#    text (i.e. at the very end of the "Problemas:" paragraph).
#
# NOTE: accented characters are built from their Unicode code points
# ([char]0x00ED = "i" with acute, [char]0x00E1 = "a" with acute) instead
# of being embedded literally, to sidestep any source-encoding pitfalls.

$iacute = [char]0x00ED
$aacute = [char]0x00E1

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the existing _GoBack bookmark (Dia 26, "...ya inicializa")
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# Step 2: find the "Problemas:" paragraph that immediately follows the
# "Diario Dia 28" heading, and append the new sentence to it.
# ---------------------------------------------------------------------
$headingText = "Diario D" + $iacute + "a 28"

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $heading = $d.Paragraphs($i).Range.Text.Trim()
    if ($heading -eq $headingText) {
        # the very next paragraph is "Problemas:"
        $target = $d.Paragraphs($i + 1)
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Problemas:' paragraph under 'Diario Dia 28'"
}

# Range covering the paragraph's text, excluding its paragraph mark.
$r = $target.Range
$body = $d.Range($r.Start, $r.End - 1)
$insertStart = $body.End

$newSentence = " Poner imagen en los botones donde " + "est" + $aacute + "n" + " las bombas."
$body.InsertAfter($newSentence)

# Give the new sentence (everything except the leading space) the same
# 18pt (sz=36) run formatting used throughout the rest of the document,
# while keeping each sentence fragment on its own run, as in the source.
$d.Range($insertStart + 1,  $insertStart + 35).Font.Size = 18   # "Poner imagen en los botones donde "
$d.Range($insertStart + 35, $insertStart + 40).Font.Size = 18   # "estan"
$d.Range($insertStart + 40, $insertStart + 51).Font.Size = 18   # " las bombas"
$d.Range($insertStart + 51, $insertStart + 52).Font.Size = 18   # "."

# ---------------------------------------------------------------------
# Step 3: re-create _GoBack collapsed right at the new end of the
# paragraph (immediately after the "."), i.e. right before the paragraph
# mark. A Range collapsed exactly on a paragraph-end boundary confuses
# Bookmarks.Add, so a one-character placeholder is inserted first to move
# the insertion point away from that boundary, and removed afterwards.
# ---------------------------------------------------------------------
$paraNow = $target.Range
$placeholderPos = $paraNow.End - 1
$d.Range($placeholderPos, $placeholderPos).InsertAfter("X")

$bmPos = $insertStart + 52
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

$d.Range($bmPos, $bmPos + 1).Delete()

Write-Output ("Updated paragraph length: " + $target.Range.Text.Length)
